$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 4.987252666666667
$ws.Range("H2").Value = 14.961758
$ws.Range("I2").Value = 0.07310203833248144
$ws.Range("J2").Value = 0.07310203833248144
$ws.Range("M2").Value = 1.442875
$ws.Range("N2").Value = 4.328625
$ws.Range("O2").Value = 0.02047893724893121
$ws.Range("P2").Value = 0.02047893724893121
$ws.Range("Q2").Value = 7.195982191416666
$ws.Range("R2").Value = 64.76383972274999
$ws.Range("S2").Value = 0.001497052055779851
$ws.Range("T2").Value = 0.001497052055779851

$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 4.987252666666667
$ws.Range("H3").Value = 14.961758
$ws.Range("I3").Value = 0.07310203833248144
$ws.Range("J3").Value = 0.07310203833248144
$ws.Range("O3").Value = 0.1473796107804731
$ws.Range("P3").Value = 0.1473796107804731
$ws.Range("Q3").Value = 51.78691851353533
$ws.Range("R3").Value = 466.082266621818
$ws.Range("S3").Value = 0.01077374995670034
$ws.Range("T3").Value = 0.01077374995670034

$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 4.987252666666667
$ws.Range("H4").Value = 14.961758
$ws.Range("I4").Value = 0.07310203833248144
$ws.Range("J4").Value = 0.07310203833248144
$ws.Range("M4").Value = 27.934719
$ws.Range("N4").Value = 83.804157
$ws.Range("O4").Value = 0.3964815784233052
$ws.Range("P4").Value = 0.3964815784233051
$ws.Range("Q4").Value = 139.317501825334
$ws.Range("R4").Value = 1253.857516428006
$ws.Range("S4").Value = 0.0289836115440232
$ws.Range("T4").Value = 0.0289836115440232

$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 4.987252666666667
$ws.Range("H5").Value = 14.961758
$ws.Range("I5").Value = 0.07310203833248144
$ws.Range("J5").Value = 0.07310203833248144
$ws.Range("M5").Value = 30.695086
$ws.Range("N5").Value = 92.085258
$ws.Range("O5").Value = 0.4356598735472906
$ws.Range("P5").Value = 0.4356598735472905
$ws.Range("Q5").Value = 153.0841495070627
$ws.Range("R5").Value = 1377.757345563564
$ws.Range("S5").Value = 0.03184762477597805
$ws.Range("T5").Value = 0.03184762477597805

$ws.Range("I6").Value = 0.7918600623966918
$ws.Range("J6").Value = 0.7918600623966917
$ws.Range("M6").Value = 1.442875
$ws.Range("N6").Value = 4.328625
$ws.Range("O6").Value = 0.02047893724893121
$ws.Range("P6").Value = 0.02047893724893121
$ws.Range("Q6").Value = 77.94872806670833
$ws.Range("R6").Value = 701.5385526003749
$ws.Range("S6").Value = 0.0162164525277566
$ws.Range("T6").Value = 0.0162164525277566

$ws.Range("I7").Value = 0.7918600623966918
$ws.Range("J7").Value = 0.7918600623966917
$ws.Range("O7").Value = 0.1473796107804731
$ws.Range("P7").Value = 0.1473796107804731
$ws.Range("Q7").Value = 560.9692077113997
$ws.Range("R7").Value = 5048.722869402597
$ws.Range("S7").Value = 0.1167040277886256
$ws.Range("T7").Value = 0.1167040277886255

$ws.Range("I8").Value = 0.7918600623966918
$ws.Range("J8").Value = 0.7918600623966917
$ws.Range("M8").Value = 27.934719
$ws.Range("N8").Value = 83.804157
$ws.Range("O8").Value = 0.3964815784233052
$ws.Range("P8").Value = 0.3964815784233051
$ws.Range("Q8").Value = 1509.122976661811
$ws.Range("R8").Value = 13582.1067899563
$ws.Range("S8").Value = 0.3139579274294173
$ws.Range("T8").Value = 0.3139579274294172

$ws.Range("I9").Value = 0.7918600623966918
$ws.Range("J9").Value = 0.7918600623966917
$ws.Range("M9").Value = 30.695086
$ws.Range("N9").Value = 92.085258
$ws.Range("O9").Value = 0.4356598735472906
$ws.Range("P9").Value = 0.4356598735472905
$ws.Range("Q9").Value = 1658.246841617067
$ws.Range("R9").Value = 14924.22157455361
$ws.Range("S9").Value = 0.3449816546508924
$ws.Range("T9").Value = 0.3449816546508923

$ws.Range("G10").Value = 9.123312666666665
$ws.Range("H10").Value = 27.369938
$ws.Range("I10").Value = 0.1337274842190096
$ws.Range("J10").Value = 0.1337274842190096
$ws.Range("M10").Value = 1.442875
$ws.Range("N10").Value = 4.328625
$ws.Range("O10").Value = 0.02047893724893121
$ws.Range("P10").Value = 0.02047893724893121
$ws.Range("Q10").Value = 13.16379976391666
$ws.Range("R10").Value = 118.47419787525
$ws.Range("S10").Value = 0.002738596757778535
$ws.Range("T10").Value = 0.002738596757778535

$ws.Range("G11").Value = 9.123312666666665
$ws.Range("H11").Value = 27.369938
$ws.Range("I11").Value = 0.1337274842190096
$ws.Range("J11").Value = 0.1337274842190096
$ws.Range("O11").Value = 0.1473796107804731
$ws.Range("P11").Value = 0.1473796107804731
$ws.Range("Q11").Value = 94.73517409695532
$ws.Range("R11").Value = 852.6165668725979
$ws.Range("S11").Value = 0.01970870457484949
$ws.Range("T11").Value = 0.01970870457484948

$ws.Range("G12").Value = 9.123312666666665
$ws.Range("H12").Value = 27.369938
$ws.Range("I12").Value = 0.1337274842190096
$ws.Range("J12").Value = 0.1337274842190096
$ws.Range("M12").Value = 27.934719
$ws.Range("N12").Value = 83.804157
$ws.Range("O12").Value = 0.3964815784233052
$ws.Range("P12").Value = 0.3964815784233051
$ws.Range("Q12").Value = 254.857175692474
$ws.Range("R12").Value = 2293.714581232266
$ws.Range("S12").Value = 0.05302048402173055
$ws.Range("T12").Value = 0.05302048402173054

$ws.Range("G13").Value = 9.123312666666665
$ws.Range("H13").Value = 27.369938
$ws.Range("I13").Value = 0.1337274842190096
$ws.Range("J13").Value = 0.1337274842190096
$ws.Range("M13").Value = 30.695086
$ws.Range("N13").Value = 92.085258
$ws.Range("O13").Value = 0.4356598735472906
$ws.Range("P13").Value = 0.4356598735472905
$ws.Range("Q13").Value = 280.0408669082226
$ws.Range("R13").Value = 2520.367802174004
$ws.Range("S13").Value = 0.058259698864651
$ws.Range("T13").Value = 0.058259698864651

$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 0.08940066666666667
$ws.Range("H14").Value = 0.268202
$ws.Range("I14").Value = 0.001310415051817319
$ws.Range("J14").Value = 0.001310415051817319
$ws.Range("M14").Value = 1.442875
$ws.Range("N14").Value = 4.328625
$ws.Range("O14").Value = 0.02047893724893121
$ws.Range("P14").Value = 0.02047893724893121
$ws.Range("Q14").Value = 0.1289939869166667
$ws.Range("R14").Value = 1.16094588225
$ws.Range("S14").Value = 0.00002683590761622182
$ws.Range("T14").Value = 0.00002683590761622181

$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 0.08940066666666667
$ws.Range("H15").Value = 0.268202
$ws.Range("I15").Value = 0.001310415051817319
$ws.Range("J15").Value = 0.001310415051817319
$ws.Range("O15").Value = 0.1473796107804731
$ws.Range("P15").Value = 0.1473796107804731
$ws.Range("Q15").Value = 0.9283237383713334
$ws.Range("R15").Value = 8.354913645342
$ws.Range("S15").Value = 0.00019312846029771
$ws.Range("T15").Value = 0.0001931284602977099

$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 0.08940066666666667
$ws.Range("H16").Value = 0.268202
$ws.Range("I16").Value = 0.001310415051817319
$ws.Range("J16").Value = 0.001310415051817319
$ws.Range("M16").Value = 27.934719
$ws.Range("N16").Value = 83.804157
$ws.Range("O16").Value = 0.3964815784233052
$ws.Range("P16").Value = 0.3964815784233051
$ws.Range("Q16").Value = 2.497382501746
$ws.Range("R16").Value = 22.476442515714
$ws.Range("S16").Value = 0.0005195554281341879
$ws.Range("T16").Value = 0.0005195554281341878

$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 0.08940066666666667
$ws.Range("H17").Value = 0.268202
$ws.Range("I17").Value = 0.001310415051817319
$ws.Range("J17").Value = 0.001310415051817319
$ws.Range("M17").Value = 30.695086
$ws.Range("N17").Value = 92.085258
$ws.Range("O17").Value = 0.4356598735472906
$ws.Range("P17").Value = 0.4356598735472905
$ws.Range("Q17").Value = 2.744161151790667
$ws.Range("R17").Value = 24.697450366116
$ws.Range("S17").Value = 0.0005708952557691995
$ws.Range("T17").Value = 0.0005708952557691993

Write-Output "done"